# Generate Report for Handback
#
# The previous handback attempt ("Ready for handoff") is replaced by a
# failure status ("Handback transform failed") across the Overview sheet
# and each per-locale sheet's Status column, and an explanatory message is
# recorded in the (previously empty) "Error Detail" column for the second
# data row of each locale sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: both the zh-cn and de-de status columns (B3, C3) show the
# status for the 11d1920b... file.
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# Per-locale sheets: column C ("Status") on row 3 holds the same status.
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Per-locale sheets: column L ("Error Detail") on row 3 gets the handback
# transform failure detail describing the generated (unexpected) handback
# file name vs. the expected handoff-derived name.
$zhcn.Range("L3").Value = "Handback file name: xwl0txjv.ey0 is different with handoff file name: 11d1920b-6b5d-4f99-8d5d-a1072b07ad82.16d7cd8ae7b4c65cec2172bfd4a57e85d753d0f2.zh-cn."
$dede.Range("L3").Value = "Handback file name: xwl0txjv.ey0 is different with handoff file name: 11d1920b-6b5d-4f99-8d5d-a1072b07ad82.16d7cd8ae7b4c65cec2172bfd4a57e85d753d0f2.de-de."
